# Update the water flow distribution table: average flow (Qavg) increased
# from 4200.4475 m3/hr to 4205.4475 m3/hr (+5 m3/hr), which shifts the
# flow-distribution table (columns B:G, rows 2:25) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (Qavg), C (Flow distribution Q), D (Difference),
# E (Amount to Storage), F (Amount from Storage), G (Running Total in Storage)
# for rows 2 through 25.
$data = @{
    2  = @(4205.447499999999, 2523.2685,             -1682.179,             0,                    1682.179,             0)
    3  = @(4205.447499999999, 2102.72375,             -2102.72375,           0,                    2102.72375,           0)
    4  = @(4205.447499999999, 1682.179,                -2523.268499999999,   0,                    2523.268499999999,   0)
    5  = @(4205.447499999999, 1261.63425,              -2943.813249999999,   0,                    2943.813249999999,   0)
    6  = @(4205.447499999999, 1261.63425,              -2943.813249999999,   0,                    2943.813249999999,   0)
    7  = @(4205.447499999999, 1682.179,                -2523.268499999999,   0,                    2523.268499999999,   0)
    8  = @(4205.447499999999, 2523.2685,                -1682.179,           0,                    1682.179,             0)
    9  = @(4205.447499999999, 2943.813249999999,        -1261.63425,         0,                    1261.63425,           0)
    10 = @(4205.447499999999, 3364.358,                 -841.0894999999991,  0,                    841.0894999999991,   0)
    11 = @(4205.447499999999, 3784.90275,                -420.5447499999996, 0,                    420.5447499999996,   0)
    12 = @(4205.447499999999, 4205.447499999999,         0,                  0,                    0,                    0)
    13 = @(4205.447499999999, 4205.447499999999,         0,                  0,                    0,                    0)
    14 = @(4205.447499999999, 4625.99225,                420.5447500000009,  420.5447500000009,    0,                    420.5447500000009)
    15 = @(4205.447499999999, 5046.536999999999,         841.0895,           841.0895,              0,                    1261.634250000001)
    16 = @(4205.447499999999, 5467.08175,                1261.634250000001,  1261.634250000001,     0,                    2523.268500000002)
    17 = @(4205.447499999999, 5887.626499999998,         1682.178999999999,  1682.178999999999,     0,                    4205.447500000001)
    18 = @(4205.447499999999, 6728.716,                  2523.268500000001,  2523.268500000001,     0,                    6728.716000000002)
    19 = @(4205.447499999999, 7149.260749999999,         2943.813249999999,  2943.813249999999,     0,                    9672.529250000001)
    20 = @(4205.447499999999, 7990.350249999999,         3784.902749999999,  3784.902749999999,     0,                    13457.432)
    21 = @(4205.447499999999, 7569.805499999999,         3364.358,           3364.358,               0,                    16821.79)
    22 = @(4205.447499999999, 6728.716,                  2523.268500000001,  2523.268500000001,     0,                    19345.0585)
    23 = @(4205.447499999999, 5046.536999999999,         841.0895,           841.0895,               0,                    20186.148)
    24 = @(4205.447499999999, 3784.90275,                -420.5447499999996, 0,                    420.5447499999996,   19765.60325)
    25 = @(4205.447499999999, 2943.813249999999,         -1261.63425,        0,                    1261.63425,           18481.969)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]   # B - Qavg
    $ws.Cells.Item($row, 3).Value = $values[1]   # C - Flow distribution Q
    $ws.Cells.Item($row, 4).Value = $values[2]   # D - Difference
    $ws.Cells.Item($row, 5).Value = $values[3]   # E - Amount to Storage
    $ws.Cells.Item($row, 6).Value = $values[4]   # F - Amount from Storage
    $ws.Cells.Item($row, 7).Value = $values[5]   # G - Running Total in Storage
}
